{"js": "// Remove the attestation paragraph (\"I personally examined the patient\n// separately and discussed the case with the resident/physician\n// assistant ...\") that appears as the first paragraph of the document\n// body. This matches the diff, which deletes that <w:p> entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target =\n  \"I personally examined the patient separately and discussed the case with the resident/physician assistant and with any services involved in a multidisciplinary fashion. I agree with the resident/physician's assistant documentation with any exceptions noted below:\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim().startsWith(target)) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the attestation paragraph (\"I personally examined the patient\n# separately and discussed the case with the resident/physician\n# assistant ...\") that appears as the first paragraph of the document\n# body. This matches the diff, which deletes that <w:p> entirely.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"I personally examined the patient separately\"\n$found = $rng.Find.Execute()\n\nif ($found) {\n    [void]$rng.Expand(4)   # wdParagraph - grow the found text to its whole paragraph (incl. mark)\n    $rng.Delete()\n}\n"}
